$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove specific account rows (identified by account number in column A).
# Delete from the bottom up so earlier row indices stay valid as we go.
$rowsToDelete = @(10, 9, 6, 2)
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}
